$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.7
$ws.Range("G3").Value = 2.86
$ws.Range("H3").Value = 2.56
$ws.Range("I3").Value = 2.7
$ws.Range("W3").Value = 1.55
$ws.Range("AA4").Value = 260
$ws.Range("AF4").Value = 10
$ws.Range("AG4").Value = 10.5
$ws.Range("AJ4").Value = 17.5
$ws.Range("AK4").Value = 18
$ws.Range("AN4").Value = 8.800000000000001
$ws.Range("AO4").Value = 170
$ws.Range("F4").Value = 1.54
$ws.Range("G4").Value = 1.59
$ws.Range("H4").Value = 6.6
$ws.Range("I4").Value = 7.6
$ws.Range("J4").Value = 4.4
$ws.Range("K4").Value = 4.8
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 4
$ws.Range("Q4").Value = 1.86
$ws.Range("R4").Value = 1.41
$ws.Range("T4").Value = 1.95
$ws.Range("U4").Value = 1.94
$ws.Range("V4").Value = 1.16
$ws.Range("W4").Value = 2.68
$ws.Range("AJ5").Value = 130
$ws.Range("F5").Value = 4.4
$ws.Range("G5").Value = 4.8
$ws.Range("H5").Value = 1.88
$ws.Range("I5").Value = 1.94
$ws.Range("J5").Value = 3.65
$ws.Range("K5").Value = 3.9
$ws.Range("V5").Value = 2.06
$ws.Range("G6").Value = 1.51
$ws.Range("J6").Value = 4.1
$ws.Range("L6").Value = 1.36
$ws.Range("AI7").Value = 75
$ws.Range("AM7").Value = 150
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 4.9
$ws.Range("K7").Value = 3.6
$ws.Range("N7").Value = 3.25
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 1.77
$ws.Range("Q7").Value = 2.2
$ws.Range("T7").Value = 1.91
$ws.Range("U7").Value = 1.96
$ws.Range("V7").Value = 1.26
$ws.Range("W7").Value = 1.98
$ws.Range("I8").Value = 8.199999999999999
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 17.5
$ws.Range("AE9").Value = 46
$ws.Range("AF9").Value = 19
$ws.Range("AH9").Value = 21
$ws.Range("AK9").Value = 29
$ws.Range("AN9").Value = 21
$ws.Range("AO9").Value = 40
$ws.Range("F9").Value = 2.18
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 3.75
$ws.Range("J9").Value = 3.4
$ws.Range("K9").Value = 4
$ws.Range("Q9").Value = 1.81
$ws.Range("T9").Value = 1.68
$ws.Range("U9").Value = 2.2
$ws.Range("V9").Value = 1.36
$ws.Range("W9").Value = 1.67
$ws.Range("Z9").Value = 30
$ws.Range("AF10").Value = 15
$ws.Range("O10").Value = 1.49
$ws.Range("S10").Value = 5.1
$ws.Range("X10").Value = 9.199999999999999
$ws.Range("Z10").Value = 19.5
$ws.Range("AK11").Value = 95
$ws.Range("AN11").Value = 140
$ws.Range("H11").Value = 1.78
$ws.Range("I11").Value = 1.95
$ws.Range("N11").Value = 3.15
$ws.Range("O11").Value = 1.37
$ws.Range("P11").Value = 1.71
$ws.Range("U11").Value = 1.85
$ws.Range("Y11").Value = 9
$ws.Range("F12").Value = 2.04
$ws.Range("G12").Value = 2.18
$ws.Range("H12").Value = 3.9
$ws.Range("I12").Value = 4.4
$ws.Range("J12").Value = 3.45
$ws.Range("K12").Value = 3.7
$ws.Range("L12").Value = 1.46
$ws.Range("S12").Value = 3.65
$ws.Range("V12").Value = 1.3
$ws.Range("W12").Value = 1.85
$ws.Range("G13").Value = 2.62
$ws.Range("I13").Value = 3.8
$ws.Range("J13").Value = 3.05
$ws.Range("M13").Value = 1.08
$ws.Range("Q13").Value = 2.22
$ws.Range("S13").Value = 4.2
$ws.Range("V13").Value = 1.35
$ws.Range("W13").Value = 1.61
$ws.Range("X13").Value = 12.5
$ws.Range("F14").Value = 3.5
$ws.Range("G14").Value = 5.3
$ws.Range("J14").Value = 3.45
$ws.Range("K14").Value = 4.4
$ws.Range("R14").Value = 1.38
$ws.Range("S14").Value = 2.46
$ws.Range("V14").Value = 1.87
$ws.Range("W14").Value = 1.23
$ws.Range("K15").Value = 3.55
$ws.Range("Q15").Value = 2.36
$ws.Range("F16").Value = 1.93
$ws.Range("G16").Value = 2.06
$ws.Range("H16").Value = 3.4
$ws.Range("K16").Value = 5.7
$ws.Range("S16").Value = 1.91
$ws.Range("V16").Value = 1.36
$ws.Range("W16").Value = 1.94
$ws.Range("AB17").Value = 25
$ws.Range("AD17").Value = 13
$ws.Range("AF17").Value = 40
$ws.Range("AK17").Value = 46
$ws.Range("AL17").Value = 46
$ws.Range("AN17").Value = 29
$ws.Range("G17").Value = 4.4
$ws.Range("H17").Value = 1.85
$ws.Range("N17").Value = 5.5
$ws.Range("O17").Value = 1.18
$ws.Range("P17").Value = 2.54
$ws.Range("Q17").Value = 1.47
$ws.Range("R17").Value = 1.62
$ws.Range("S17").Value = 2.16
$ws.Range("T17").Value = 1.5
$ws.Range("U17").Value = 2.5
$ws.Range("W17").Value = 1.29
$ws.Range("X17").Value = 32
$ws.Range("Y17").Value = 16
$ws.Range("Z17").Value = 17.5
$ws.Range("J18").Value = 3.1
$ws.Range("L18").Value = 1.46
$ws.Range("W18").Value = 1.98
$ws.Range("G19").Value = 2.62
$ws.Range("J19").Value = 3.05
$ws.Range("R19").Value = 1.21
$ws.Range("T19").Value = 2.06
$ws.Range("F20").Value = 2.7
$ws.Range("I20").Value = 2.92
$ws.Range("T20").Value = 1.85
$ws.Range("U20").Value = 1.99
$ws.Range("V20").Value = 1.51
$ws.Range("AF21").Value = 17.5
$ws.Range("AL21").Value = 60
$ws.Range("G21").Value = 2.96
$ws.Range("N21").Value = 3
$ws.Range("T21").Value = 2.04
$ws.Range("X21").Value = 9.4
$ws.Range("AA22").Value = 80
$ws.Range("AD22").Value = 16.5
$ws.Range("AG22").Value = 10
$ws.Range("T22").Value = 1.53
$ws.Range("U22").Value = 2.8
$ws.Range("AH23").Value = 19
$ws.Range("G23").Value = 2.34
$ws.Range("Q23").Value = 2.22
$ws.Range("W23").Value = 1.74
$ws.Range("F24").Value = 2.88
$ws.Range("G24").Value = 3
$ws.Range("I24").Value = 2.72
$ws.Range("L24").Value = 1.4
$ws.Range("Q24").Value = 1.96
$ws.Range("R24").Value = 1.37
$ws.Range("V24").Value = 1.58
$ws.Range("AA25").Value = 80
$ws.Range("AB25").Value = 8.4
$ws.Range("AC25").Value = 7.4
$ws.Range("AD25").Value = 1000
$ws.Range("AE25").Value = 60
$ws.Range("AF25").Value = 1000
$ws.Range("AG25").Value = 13.5
$ws.Range("AH25").Value = 980
$ws.Range("AJ25").Value = 40
$ws.Range("AN25").Value = 980
$ws.Range("F25").Value = 2.26
$ws.Range("G25").Value = 2.44
$ws.Range("I25").Value = 4
$ws.Range("L25").Value = 1.5
$ws.Range("M25").Value = 1.1
$ws.Range("P25").Value = 1.67
$ws.Range("Q25").Value = 2.28
$ws.Range("S25").Value = 4.4
$ws.Range("T25").Value = 1.95
$ws.Range("U25").Value = 1.92
$ws.Range("W25").Value = 1.69
$ws.Range("X25").Value = 11
$ws.Range("Y25").Value = 1000
$ws.Range("Z25").Value = 980
